$wb = $excel.ActiveWorkbook

# --- Sheet 1: AddRecipe ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "Pizza"
$ws1.Range("B2").Value = "Chicken Fajita"
$ws1.Range("C2").Value = "20 Mins"

# --- Sheet 2: EditRecipe ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Pizza"
$ws2.Range("B2").Value = "Chicken Fajita"
$ws2.Range("C2").Value = "18 mins"
$ws2.Range("D2").Value = "Cheese"
$ws2.Range("E2").Value = "250 grams"

# --- Sheet 3: DeleteRecipe ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Pizza"

# --- Sheet 4: AddCategory ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = "Pizza"
$ws4.Range("B2").Value = "Entrée"

# --- Page setup (sheet1 & sheet4 switch to portrait orientation) ---
$ws1.PageSetup.Orientation = 1
$ws4.PageSetup.Orientation = 1

# --- Restore selections / active sheet, in the order the author left them ---
$ws1.Range("A2:B2").Select() | Out-Null
$ws2.Range("F11").Select() | Out-Null
$ws3.Range("B3").Select() | Out-Null
$ws4.Range("B2").Select() | Out-Null
